# Applies the crypto price/volume refresh described by the commit diff.
# Column D values that parse as plain numbers are written with a leading
# apostrophe so Excel stores them as text (matching the original inlineStr/
# shared-string cells) instead of silently coercing them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '90.811.52'
$ws.Cells.Item(2, 5).Value = '  +1.23%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.161.01'
$ws.Cells.Item(3, 5).Value = '  +3.75%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.23%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''215.49'
$ws.Cells.Item(5, 5).Value = '  +1.83%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''627.77'
$ws.Cells.Item(6, 5).Value = '  +2.40%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +28.88%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +3.22%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.00%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '3.156.41'
$ws.Cells.Item(10, 5).Value = '  +3.69%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.752'
$ws.Cells.Item(11, 5).Value = '  +11.26%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +8.51%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'ShibaInu'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(13, 4).Value = '''0.0000246'
$ws.Cells.Item(13, 5).Value = '  +2.81%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'Toncoin'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(14, 4).Value = '''5.69'
$ws.Cells.Item(14, 5).Value = '  +5.77%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''34.96'
$ws.Cells.Item(15, 5).Value = '  +8.22%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '90.541.30'
$ws.Cells.Item(16, 5).Value = '  +1.08%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '3.740.08'
$ws.Cells.Item(17, 5).Value = '  +4.16%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.144.72'
$ws.Cells.Item(18, 5).Value = '  +3.25%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +11.24%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''14.32'
$ws.Cells.Item(20, 5).Value = '  +6.77%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  -2.84%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''466.27'
$ws.Cells.Item(22, 5).Value = '  +9.81%  '

# Row 23
$ws.Cells.Item(23, 2).Value = 'Polkadot'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(23, 4).Value = '''5.56'
$ws.Cells.Item(23, 5).Value = '  +10.70%  '

# Row 24
$ws.Cells.Item(24, 2).Value = 'Uniswap'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(24, 4).Value = '''9.08'
$ws.Cells.Item(24, 5).Value = '  +10.06%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''5.88'
$ws.Cells.Item(25, 5).Value = '  +9.23%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''93.63'
$ws.Cells.Item(26, 5).Value = '  +11.62%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''12.18'
$ws.Cells.Item(27, 5).Value = '  +4.61%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '3.316.55'
$ws.Cells.Item(28, 5).Value = '  +3.59%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.09%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +2.83%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -0.56%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''9.14'
$ws.Cells.Item(32, 5).Value = '  +10.25%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''27.10'
$ws.Cells.Item(33, 5).Value = '  +18.52%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''522.01'
$ws.Cells.Item(34, 5).Value = '  +3.65%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''0.182'
$ws.Cells.Item(35, 5).Value = '  +33.34%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -1.99%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +7.36%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''6.92'
$ws.Cells.Item(38, 5).Value = '  +3.92%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'Kaspa'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(39, 4).Value = '''0.142'
$ws.Cells.Item(39, 5).Value = '  +7.87%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'Fetch.AI'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(40, 4).Value = '''1.31'
$ws.Cells.Item(40, 5).Value = '  +5.36%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(41, 4).Value = '''22.21'
$ws.Cells.Item(41, 5).Value = '  -0.13%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Hedera'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(42, 4).Value = '''0.0862'
$ws.Cells.Item(42, 5).Value = '  +25.05%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +0.01%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''0.416'
$ws.Cells.Item(44, 5).Value = '  +14.58%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +7.74%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +0.00%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''150.26'
$ws.Cells.Item(47, 5).Value = '  +2.40%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''45.24'
$ws.Cells.Item(48, 5).Value = '  +4.47%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Filecoin'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(49, 4).Value = '''4.53'
$ws.Cells.Item(49, 5).Value = '  +7.51%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'ImmutableX'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(50, 4).Value = '''1.35'
$ws.Cells.Item(50, 5).Value = '  +11.05%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''0.678'
$ws.Cells.Item(51, 5).Value = '  +15.29%  '
